# Remove the stale "Ver no Jupiter..." and copyright footer paragraphs
# (plus the blank paragraph separating them from the preceding
# "LOQ4031: Química Geral I (Requisito fraco)" requirement line), while
# leaving everything else -- including the trailing blank / page-break
# paragraphs -- untouched.

$d = $word.ActiveDocument

# Locate the "Ver no Jupiter..." paragraph; the blank paragraph right
# before it (the one that currently follows the LOQ4031 line) is where
# the deletion should start.
$jupiterFind = $d.Content
$jupiterFind.Find.Execute("Ver no Jupiter Salvar em pdf Salvar em docx", `
    $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$beforeJupiter = $d.Range(0, $jupiterFind.Start)
$blankIndex = $beforeJupiter.Paragraphs.Count
$startPos = $d.Paragraphs.Item($blankIndex).Range.Start

# Locate the copyright paragraph; deletion ends at the end of it.
$copyrightFind = $d.Content
$copyrightFind.Find.Execute("Contact: luizeleno@usp.br", `
    $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$beforeCopyright = $d.Range(0, $copyrightFind.Start)
$copyrightIndex = $beforeCopyright.Paragraphs.Count
$endPos = $d.Paragraphs.Item($copyrightIndex).Range.End

# Delete the blank paragraph + "Ver no Jupiter..." paragraph +
# copyright paragraph in one shot.
$d.Range($startPos, $endPos).Delete()
